$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(1,1).Range.Text = "42×72=3024"
$t.Cell(1,2).Range.Text = "24×88=2112"
$t.Cell(1,3).Range.Text = "30×26=780"
$t.Cell(1,4).Range.Text = "81×23=1863"
$t.Cell(1,5).Range.Text = "93×22=2046"
$t.Cell(2,1).Range.Text = "43×29=1247"
$t.Cell(2,2).Range.Text = "77×94=7238"
$t.Cell(2,3).Range.Text = "41×74=3034"
$t.Cell(2,4).Range.Text = "90×11=990"
$t.Cell(2,5).Range.Text = "25×10=250"
$t.Cell(3,1).Range.Text = "45×99=4455"
$t.Cell(3,2).Range.Text = "53×40=2120"
$t.Cell(3,3).Range.Text = "58×54=3132"
$t.Cell(3,4).Range.Text = "97×38=3686"
$t.Cell(3,5).Range.Text = "41×36=1476"
$t.Cell(4,1).Range.Text = "15×54=810"
$t.Cell(4,2).Range.Text = "15×85=1275"
$t.Cell(4,3).Range.Text = "62×20=1240"
$t.Cell(4,4).Range.Text = "84×51=4284"
$t.Cell(4,5).Range.Text = "11×59=649"
$t.Cell(5,1).Range.Text = "53×80=4240"
$t.Cell(5,2).Range.Text = "45×84=3780"
$t.Cell(5,3).Range.Text = "34×18=612"
$t.Cell(5,4).Range.Text = "29×94=2726"
$t.Cell(5,5).Range.Text = "83×13=1079"
$t.Cell(6,1).Range.Text = "52×33=1716"
$t.Cell(6,2).Range.Text = "76×79=6004"
$t.Cell(6,3).Range.Text = "59×41=2419"
$t.Cell(6,4).Range.Text = "77×20=1540"
$t.Cell(6,5).Range.Text = "79×46=3634"
$t.Cell(7,1).Range.Text = "57×40=2280"
$t.Cell(7,2).Range.Text = "86×70=6020"
$t.Cell(7,3).Range.Text = "64×33=2112"
$t.Cell(7,4).Range.Text = "16×75=1200"
$t.Cell(7,5).Range.Text = "96×34=3264"
$t.Cell(8,1).Range.Text = "36×90=3240"
$t.Cell(8,2).Range.Text = "34×59=2006"
$t.Cell(8,3).Range.Text = "16×55=880"
$t.Cell(8,4).Range.Text = "24×70=1680"
$t.Cell(8,5).Range.Text = "86×69=5934"
$t.Cell(9,1).Range.Text = "92×14=1288"
$t.Cell(9,2).Range.Text = "54×50=2700"
$t.Cell(9,3).Range.Text = "32×100=3200"
$t.Cell(9,4).Range.Text = "30×99=2970"
$t.Cell(9,5).Range.Text = "73×24=1752"
$t.Cell(10,1).Range.Text = "54×61=3294"
$t.Cell(10,2).Range.Text = "26×100=2600"
$t.Cell(10,3).Range.Text = "89×40=3560"
$t.Cell(10,4).Range.Text = "15×34=510"
$t.Cell(10,5).Range.Text = "41×78=3198"
$t.Cell(11,1).Range.Text = "39×92=3588"
$t.Cell(11,2).Range.Text = "39×42=1638"
$t.Cell(11,3).Range.Text = "77×24=1848"
$t.Cell(11,4).Range.Text = "83×60=4980"
$t.Cell(11,5).Range.Text = "56×97=5432"
$t.Cell(12,1).Range.Text = "37×26=962"
$t.Cell(12,2).Range.Text = "59×79=4661"
$t.Cell(12,3).Range.Text = "46×39=1794"
$t.Cell(12,4).Range.Text = "47×86=4042"
$t.Cell(12,5).Range.Text = "23×57=1311"
$t.Cell(13,1).Range.Text = "12×52=624"
$t.Cell(13,2).Range.Text = "51×42=2142"
$t.Cell(13,3).Range.Text = "30×33=990"
$t.Cell(13,4).Range.Text = "52×10=520"
$t.Cell(13,5).Range.Text = "68×42=2856"
$t.Cell(14,1).Range.Text = "65×10=650"
$t.Cell(14,2).Range.Text = "43×24=1032"
$t.Cell(14,3).Range.Text = "14×30=420"
$t.Cell(14,4).Range.Text = "73×86=6278"
$t.Cell(14,5).Range.Text = "97×70=6790"
$t.Cell(15,1).Range.Text = "88×72=6336"
$t.Cell(15,2).Range.Text = "33×97=3201"
$t.Cell(15,3).Range.Text = "73×23=1679"
$t.Cell(15,4).Range.Text = "66×11=726"
$t.Cell(15,5).Range.Text = "12×91=1092"
$t.Cell(16,1).Range.Text = "38×21=798"
$t.Cell(16,2).Range.Text = "80×39=3120"
$t.Cell(16,3).Range.Text = "51×36=1836"
$t.Cell(16,4).Range.Text = "82×45=3690"
$t.Cell(16,5).Range.Text = "40×100=4000"
$t.Cell(17,1).Range.Text = "92×70=6440"
$t.Cell(17,2).Range.Text = "53×67=3551"
$t.Cell(17,3).Range.Text = "50×12=600"
$t.Cell(17,4).Range.Text = "75×82=6150"
$t.Cell(17,5).Range.Text = "70×97=6790"
$t.Cell(18,1).Range.Text = "67×61=4087"
$t.Cell(18,2).Range.Text = "11×75=825"
$t.Cell(18,3).Range.Text = "98×85=8330"
$t.Cell(18,4).Range.Text = "58×91=5278"
$t.Cell(18,5).Range.Text = "38×50=1900"
$t.Cell(19,1).Range.Text = "21×83=1743"
$t.Cell(19,2).Range.Text = "73×80=5840"
$t.Cell(19,3).Range.Text = "45×97=4365"
$t.Cell(19,4).Range.Text = "80×89=7120"
$t.Cell(19,5).Range.Text = "20×30=600"
$t.Cell(20,1).Range.Text = "91×31=2821"
$t.Cell(20,2).Range.Text = "66×18=1188"
$t.Cell(20,3).Range.Text = "86×85=7310"
$t.Cell(20,4).Range.Text = "25×93=2325"
$t.Cell(20,5).Range.Text = "42×15=630"
